$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the tab strip (after the last existing
# sheet) and activate it -- this both re-points tabSelected/activeTab at the
# new sheet and removes it from the previously active "Box Mix array" sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Number Fmt"
$ws.Activate()

# Row 1
$ws.Range("A1").Value = 5
$ws.Range("B1").Value = 50000
$ws.Range("C1").Value = 0.005
$ws.Range("D1").Value = 0.0000005
$ws.Range("E1").Value = 500000000
$ws.Range("F1").Value = 500000000.5

# Row 2
$ws.Range("A2").Value = -5
$ws.Range("B2").Value = -50000
$ws.Range("C2").Value = -0.005
$ws.Range("D2").Value = -0.0000005
$ws.Range("E2").Value = -500000000
$ws.Range("F2").Value = -500000000.5

# Scientific-notation number format (built-in numFmtId 11 => "0.00E+00")
# applied to everything but A1:A2 and B2 (matching the original author's
# worksheet, which left those cells on General format).
$ws.Range("B1").NumberFormat = "0.00E+00"
$ws.Range("C1:C2").NumberFormat = "0.00E+00"
$ws.Range("D1:D2").NumberFormat = "0.00E+00"
$ws.Range("E1:E2").NumberFormat = "0.00E+00"
$ws.Range("F1:F2").NumberFormat = "0.00E+00"

# Columns E:F were widened/best-fit to show the longer scientific values.
$ws.Columns("E:F").ColumnWidth = 9.375

# Page setup (paper size / orientation) mirroring the other sheets.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Mirror the row-by-row entry flow: after typing the data the user is left
# with the selection on A3 (one past the last populated row).
$null = $ws.Range("A3").Select()
